$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "Fork-Join" timing column (H) and its SpeedUp column (I),
#     mirroring the existing "Thread Java" (F) / "SpeedUp" (G) pair ---
$ws.Range("H1").Value = "Fork-Join"
$ws.Range("I1").Value = "SpeedUp"

# --- Row 2 (PI) is marked "NA" for both new columns, like it already is for F/G ---
$ws.Range("H2").Value = "NA"
$ws.Range("H2").NumberFormat = "0.00"
$ws.Range("H2").HorizontalAlignment = -4152
$ws.Range("I2").Value = "NA"
$ws.Range("I2").NumberFormat = "0.00"
$ws.Range("I2").HorizontalAlignment = -4152

# --- Row 3 (SUM): Fork-Join timing + SpeedUp formula ---
$ws.Range("H3").Value = 85
$ws.Range("H3").NumberFormat = "0.00"
$ws.Range("I3").Formula = "=D3/H3"
$ws.Range("I3").NumberFormat = "0.00"

# --- Row 4 (2VEC) ---
$ws.Range("H4").Value = 122
$ws.Range("H4").NumberFormat = "0.00"
$ws.Range("I4").Formula = "=D4/H4"
$ws.Range("I4").NumberFormat = "0.00"

# --- Row 5 (COUNTING) ---
$ws.Range("H5").Value = 102.2
$ws.Range("H5").NumberFormat = "0.00"
$ws.Range("I5").Formula = "=D5/H5"
$ws.Range("I5").NumberFormat = "0.00"

# stray formatted (underlined) blank cell a couple of columns over
$ws.Range("K5").Font.Underline = 2
$ws.Range("K5").Font.Name = "Arial"
$ws.Range("K5").Font.Size = 10

# --- Rows 6-8 (JULIA, BLUR, FACTORIAL): columns not filled in yet for Fork-Join,
#     but formatted the same as the rest of the H/I columns ---
$ws.Range("H6").NumberFormat = "0.00"
$ws.Range("I6").NumberFormat = "0.00"
$ws.Range("H7").NumberFormat = "0.00"
$ws.Range("I7").NumberFormat = "0.00"
$ws.Range("H8").NumberFormat = "0.00"
$ws.Range("I8").NumberFormat = "0.00"

$ws.Range("I6").Select()
